$wb = $excel.ActiveWorkbook

$wsProcs = $wb.Worksheets.Item("ELC_IMP_EXP_ProcsR")
$wsTechs = $wb.Worksheets.Item("ELC_IMP_EXP_TechsR")

# ----------------------------------------------------------------------
# ELC_IMP_EXP_ProcsR: fix the IMP row (was row 5, too close to the
# insert-table placeholder) -- drop the leading '*' from the tech name
# and correct the unit column, then pull the row up to row 4.
# ----------------------------------------------------------------------
$wsProcs.Range("B5").Value = "IMP"
$wsProcs.Range("F5").Value = "Pja"
$wsProcs.Rows.Item(4).Delete()

# ----------------------------------------------------------------------
# ELC_IMP_EXP_TechsR: add the VAROM column, fix the IMPELC-DKW technology
# (drop leading '*', correct capacity/efficiency), then pull the row up
# to row 4 as well.
# ----------------------------------------------------------------------
$wsTechs.Range("H2").Copy()
$wsTechs.Range("J2").PasteSpecial(-4122)
$wsTechs.Range("J2").Borders.LineStyle = -4142
$wsTechs.Range("J2").Value = "VAROM"

$wsTechs.Range("I3").Copy()
$wsTechs.Range("J3").PasteSpecial(-4122)
$wsTechs.Range("J3").Value = 0

$wsTechs.Range("B5").Value = "IMPELC-DKW"
$wsTechs.Range("F5").Value = 100
$wsTechs.Range("H5").Value = 1

$wsTechs.Range("I5").Copy()
$wsTechs.Range("J5").PasteSpecial(-4122)
$wsTechs.Range("J5").Value = 1

$wsTechs.Rows.Item(4).Delete()

# ----------------------------------------------------------------------
# Reposition the comment text box on ELC_IMP_EXP_TechsR so it no longer
# overlaps the (now shifted) insert-table area.
# ----------------------------------------------------------------------
$shp = $wsTechs.Shapes.Item(1)
$shp.Left = 11049000 / 12700
$shp.Top = 160020 / 12700
$shp.Width = 2415540 / 12700
$shp.Height = 2141220 / 12700

# ----------------------------------------------------------------------
# Selections / active sheet: ELC_IMP_EXP_TechsR becomes the active tab.
# ----------------------------------------------------------------------
$wsProcs.Activate()
$wsProcs.Range("E5").Select()

$wsTechs.Activate()
$wsTechs.Range("K5").Select()
